# Apply cryptos list update (values sourced from the target OOXML diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.674.56"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "3.477.27"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("D7").Value = "3.476.41"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "4.062.48"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "3.481.83"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "63.695.54"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "3.614.38"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000111"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").Value = "3.477.02"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0795"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "2.424.04"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("E51").Value = "  -1.02%  "
